# Remove the "Dose multiplier interval" column (column G).
# This shifts Dose vol. (was H) into G and Force delay (was I) into H,
# removes the now-unused shared string, shrinks the used range from
# A1:AF12 to A1:AE12, and updates the column width definitions
# accordingly - mirroring an in-Excel "Delete Column" on G.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").EntireColumn.Delete()

# Reflect the selection left on column G (now "Dose vol.") after the
# delete, matching a full-column selection as shown by the author.
$ws.Columns("G:G").Select()
